# Generate Report for Handoff
# Updates the localization-status report:
#  - Refreshes the "Latest HO Xliff Generate Date" timestamp on the Overview sheet
#  - Refreshes the "Latest Handoff Datetime" timestamp on the zh-cn sheet
#  - Marks Priority = "ht" for the files that just went through handoff,
#    on both the zh-cn and de-de sheets

$wb = $excel.ActiveWorkbook

$rows = @(7, 8, 10, 11, 13, 14)

# --- Overview sheet: Latest HO Xliff Generate Date (column G) ---
$overview = $wb.Worksheets.Item("Overview")
foreach ($r in $rows) {
    $overview.Range("G$r").Value = "2016-09-05 10:27:39"
}

# --- zh-cn sheet: Latest Handoff Datetime (column H) + Priority (column E) ---
$zhcn = $wb.Worksheets.Item("zh-cn")
foreach ($r in $rows) {
    $zhcn.Range("H$r").Value = "2016-09-05 10:27:33"
    $zhcn.Range("E$r").Value = "ht"
}

# --- de-de sheet: Priority (column E) ---
# Column H on this sheet shares the same underlying text as Overview's
# column G ("Latest HO Xliff Generate Date"), so it must be refreshed too.
$dede = $wb.Worksheets.Item("de-de")
foreach ($r in $rows) {
    $dede.Range("H$r").Value = "2016-09-05 10:27:39"
    $dede.Range("E$r").Value = "ht"
}
